$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the RAD test data timestamps in column B (rows 2-4)
$ws.Range("B2").Value = "Tue Feb 11 20:22:49 EST 2025"
$ws.Range("B3").Value = "Tue Feb 11 20:23:03 EST 2025"
$ws.Range("B4").Value = "Tue Feb 11 20:23:16 EST 2025"
